# adding exit columns, voiding in exports
#
# Adds three new trailing header columns (exit_status, exit_date,
# exit_reason) to the DREAMS export template's single header row, and
# widens the last pre-existing column (GU, index 203) slightly so the
# new "exit_*" columns have room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells appended right after the previous last column (GU = 203).
$ws.Cells.Item(1, 204).Value = "exit_status"
$ws.Cells.Item(1, 205).Value = "exit_date"
$ws.Cells.Item(1, 206).Value = "exit_reason"

# Column 203 (GU) picks up an explicit (non bestFit) width in the target
# workbook; land on the closest value the host's pixel-quantized
# ColumnWidth model can represent.
$ws.Columns.Item(203).ColumnWidth = 11.4
